# New accelerometer readings (x, y, z) were captured and prepended to the
# top of the dataset. This pushes the previously-recorded rows down by four
# rows, and the four oldest rows that fall past the bottom of the sheet are
# removed so the sheet keeps the same 20-row window of samples (A1:C21
# instead of the original A1:C22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final x,y,z values for data rows 2-21 (header stays in row 1):
#   rows 2-5   -> four brand new accelerometer samples
#   rows 6-21  -> the samples that used to be in rows 2-17
$allData = @(
    @(-1.03581714630127,   1.381664276123047,   -0.1787742376327515),
    @(-1.100839495658875,  1.406517148017883,   -0.2175595723092557),
    @(-1.232075214385986,  1.368059515953064,   -0.2088889628648757),
    @(-1.198645412921905,  1.359035015106201,   -0.2234921492636203),
    @(-1.162086248397827,  1.317891120910645,   0.1167446374893195),
    @(-1.455766379833223,  1.126043200492858,   0.7863338142633457),
    @(-1.653480172157284,  0.8844107389450073,  1.495913922786714),
    @(-1.116286456584938,  0.6627160906791623,  2.260028153657919),
    @(-3.20666265487671,   -0.9849638938903851, 3.828832626342773),
    @(-3.492754817008973,  -1.960709273815156,  3.70500636100769),
    @(-2.38550305366516,   0.4144415855407741,  -0.4261573851108604),
    @(-2.568133831024171,  3.730020523071294,   -0.8093817904591537),
    @(-1.976609468460079,  6.545797109603887,   2.459241539239894),
    @(-1.197644114494325,  5.058232277631741,   2.450196892023069),
    @(-1.255056142807007,  0.5559926331043243,  -1.602012172341346),
    @(-0.5990372896194439, 1.264416024088864,   -0.3774302378296863),
    @(-0.1129188537597673, 2.17881894111634,    -0.6807380914687996),
    @(-0.4983874559402485, 3.747065991163262,   1.904177859425558),
    @(-0.2239453792572,    0.9145344123244095,  -1.341536760330224),
    @(0.0572257041931156,  1.018226306885483,   0.3575173318386149)
)

for ($i = 0; $i -lt $allData.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $allData[$i][0]
    $ws.Cells.Item($row, 2).Value = $allData[$i][1]
    $ws.Cells.Item($row, 3).Value = $allData[$i][2]
}

# Row 22 held the oldest sample and is no longer part of the window; clear
# it so the sheet's used range shrinks back down to A1:C21.
$ws.Range("A22:C22").ClearContents() | Out-Null
